$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.538.52"
$ws.Range("E2").Value = "  +4.00%  "

# Row 3
$ws.Range("D3").Value = "2.697.74"
$ws.Range("E3").Value = "  +3.81%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.80%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.62%  "

# Row 7
$ws.Range("E7").Value = "  -0.35%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.605"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.17%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.67"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.67%  "

# Row 10
$ws.Range("E10").Value = "  +6.46%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.385"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.49%  "

# Row 12
$ws.Range("E12").Value = "  +1.21%  "

# Row 13
$ws.Range("D13").Value = "3.163.08"
$ws.Range("E13").Value = "  +3.34%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.90%  "

# Row 15
$ws.Range("D15").Value = "62.410.60"
$ws.Range("E15").Value = "  +3.79%  "

# Row 16
$ws.Range("E16").Value = "  +6.34%  "

# Row 17
$ws.Range("D17").Value = "2.689.41"
$ws.Range("E17").Value = "  +3.29%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.52%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.99%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "362.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.83%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.93%  "

# Row 22
$ws.Range("E22").Value = "  +0.31%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.532"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.29%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.69%  "

# Row 25
$ws.Range("E25").Value = "  +3.90%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.68%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.08%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.45%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0843"
$ws.Range("E29").Value = "  +6.26%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +11.14%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "170.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.25%  "

# Row 32
$ws.Range("E32").Value = "  -0.17%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.48"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.55%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +18.60%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.74"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.12%  "

# Row 36
$ws.Range("E36").Value = "  +8.33%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.70%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +20.00%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "353.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +13.55%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.52%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.93%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +12.88%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.62%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0588"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.95%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.62%  "

# Row 46
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.638"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.62%  "

# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "137.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.19%  "

# Row 48
$ws.Range("E48").Value = "  +6.42%  "

# Row 49
$ws.Range("E49").Value = "  +1.47%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.995"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.46%  "

# Row 51
$ws.Range("D51").Value = "2.130.37"
$ws.Range("E51").Value = "  +5.99%  "
